# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (E) and "Correspond Handback DateTime" (H)
# columns for the 968e5e94-... row (row 3) on both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-24 06:53:48"
$wsZhCn.Range("H3").Value = "2016-03-24 06:54:21"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-24 06:53:53"
$wsDeDe.Range("H3").Value = "2016-03-24 06:54:27"
